# Applies the textual content changes described by the commit diff.
# (Cosmetic-only changes in the diff -- run splitting/merging introduced by
#  proofErr tags, lastRenderedPageBreak repositioning, empty paragraph-mark
#  rPr removal, and the customXml GUID -- are artifacts of Word re-saving
#  the file and don't change the document's visible text, so they are not
#  reproduced here.)

$d = $word.ActiveDocument

# 1. Insert a parenthetical clarification after "...欄所示" and before the
#    closing paren of the preceding "(詳見附表二編號..." citation.
$d.Content.Find.Execute(
    "至66「約定紅利」欄所示)，確實存有相當比例之差異",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "至66「約定紅利」欄所示(泛指判決書編號))，確實存有相當比例之差異",
    2) | Out-Null

# 2. Drop the "實" before "損害金融秩序甚鉅" and add a closing remark before
#    the sentence's full stop.
$d.Content.Find.Execute(
    "，實損害金融秩序甚鉅。",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "，損害金融秩序甚鉅，實為二次吸金主謀。",
    2) | Out-Null

# 3. Flesh out the abbreviated court reference with the full case number.
$d.Content.Find.Execute(
    "經臺南地方法院刑事一審判決其違反銀行法",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "經臺灣臺南地方法院 108 年度金重訴字第3號刑事判決其違反銀行法",
    2) | Out-Null

# 4. Drop the leading "經查，" before citing the Supreme Court holding.
$d.Content.Find.Execute(
    "經查，最高法院98年度台上字第2157號判決之判決",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "最高法院98年度台上字第2157號判決之判決",
    2) | Out-Null

# 5. Rephrase the claim basis sentence.
$d.Content.Find.Execute(
    "之侵權行為請求上訴人及視同上訴人連帶賠償",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "之侵權行為主張上訴人針對侵權行為損害賠償",
    2) | Out-Null

# 6. Remove "連帶" (joint-and-several) from the refusal-to-pay clause.
$d.Content.Find.Execute(
    "也不願連帶賠償被上訴人損失",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "也不願賠償被上訴人損失",
    2) | Out-Null
